$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template row 422 (A:L) always has a value in every column A-L, used as the
# formatting template for the new rows so the new cells inherit the same
# cell styles (s="2" date, s="3" text/number, s="4" percent) as the rest of the sheet.
$templateAL = $ws.Range("A422:L422")
# Row 422 has "M" filled (Red branch); row 421 has "N" filled (Black branch).
# Copying from whichever one matches the branch for each new row gives the new
# M/N cell the correct style without ever materializing the unused column.
$templateM = $ws.Range("M422")
$templateN = $ws.Range("N421")

# Row 423
$templateAL.Copy($ws.Range("A423:L423"))
$templateM.Copy($ws.Range("M423"))
$ws.Range("A423").Value = 45200.93687467593
$ws.Range("B423").Value = "sumin102573@naver.com"
$ws.Range("C423").Value = "경영학과"
$ws.Range("D423").Value = 20212922
$ws.Range("E423").Value = "김수민"
$ws.Range("F423").Value = "민주 문자"
$ws.Range("G423").Value = "한글"
$ws.Range("H423").Value = "하나도 없다"
$ws.Range("I423").Value = 0.9
$ws.Range("J423").Value = "미국"
$ws.Range("K423").Value = "건강이 좋지 않다"
$ws.Range("L423").Value = "Red"
$ws.Range("M423").Value = "휴우, 그래도 반이나 남았네."

# Row 424
$templateAL.Copy($ws.Range("A424:L424"))
$templateM.Copy($ws.Range("M424"))
$ws.Range("A424").Value = 45200.93968752315
$ws.Range("B424").Value = "limmh96@gmail.com"
$ws.Range("C424").Value = "광고홍보학과"
$ws.Range("D424").Value = 20232632
$ws.Range("E424").Value = "임민호"
$ws.Range("F424").Value = "엘리트 문자"
$ws.Range("G424").Value = "한글"
$ws.Range("H424").Value = "2개"
$ws.Range("I424").Value = 0.8
$ws.Range("J424").Value = "미국"
$ws.Range("K424").Value = "시간당 중위 임금이 60% 낮다"
$ws.Range("L424").Value = "Red"
$ws.Range("M424").Value = "휴우, 그래도 반이나 남았네."

# Row 425
$templateAL.Copy($ws.Range("A425:L425"))
$templateN.Copy($ws.Range("N425"))
$ws.Range("A425").Value = 45200.94147167824
$ws.Range("B425").Value = "scw0922@naver.com"
$ws.Range("C425").Value = "간호학과"
$ws.Range("D425").Value = 20236256
$ws.Range("E425").Value = "신채원"
$ws.Range("F425").Value = "민주 문자"
$ws.Range("G425").Value = "한글"
$ws.Range("H425").Value = "하나도 없다"
$ws.Range("I425").Value = 0.5
$ws.Range("J425").Value = "미국"
$ws.Range("K425").Value = "남들을 덜 신뢰한다"
$ws.Range("L425").Value = "Black"
$ws.Range("N425").Value = "헐, 반 밖에 안 남았네."

# Row 426
$templateAL.Copy($ws.Range("A426:L426"))
$templateM.Copy($ws.Range("M426"))
$ws.Range("A426").Value = 45200.94627571759
$ws.Range("B426").Value = "harin3040@naver.com"
$ws.Range("C426").Value = "심리학과"
$ws.Range("D426").Value = 20232113
$ws.Range("E426").Value = "김현진"
$ws.Range("F426").Value = "엘리트 문자"
$ws.Range("G426").Value = "한글"
$ws.Range("H426").Value = "하나도 없다"
$ws.Range("I426").Value = 0.8
$ws.Range("J426").Value = "이탈리아"
$ws.Range("K426").Value = "시간당 중위 임금이 60% 낮다"
$ws.Range("L426").Value = "Red"
$ws.Range("M426").Value = "헐, 반 밖에 안 남았네."

# Row 427
$templateAL.Copy($ws.Range("A427:L427"))
$templateM.Copy($ws.Range("M427"))
$ws.Range("A427").Value = 45200.947467094906
$ws.Range("B427").Value = "shanesun0923@gmail.com"
$ws.Range("C427").Value = "간호학과"
$ws.Range("D427").Value = 20236253
$ws.Range("E427").Value = "선세인"
$ws.Range("F427").Value = "민주 문자"
$ws.Range("G427").Value = "한글"
$ws.Range("H427").Value = "3개"
$ws.Range("I427").Value = 0.9
$ws.Range("J427").Value = "대한민국"
$ws.Range("K427").Value = "사회활동이나 자원활동에 덜 참여한다"
$ws.Range("L427").Value = "Red"
$ws.Range("M427").Value = "모름/기타"

# Row 428
$templateAL.Copy($ws.Range("A428:L428"))
$templateM.Copy($ws.Range("M428"))
$ws.Range("A428").Value = 45200.95396979166
$ws.Range("B428").Value = "sung93716@gmail.com"
$ws.Range("C428").Value = "데이터사이언스학부"
$ws.Range("D428").Value = 20233261
$ws.Range("E428").Value = "한예림"
$ws.Range("F428").Value = "민주 문자"
$ws.Range("G428").Value = "한글"
$ws.Range("H428").Value = "2개"
$ws.Range("I428").Value = 0.1
$ws.Range("J428").Value = "이탈리아"
$ws.Range("K428").Value = "2배 정도 실직할 가능성이 높다"
$ws.Range("L428").Value = "Red"
$ws.Range("M428").Value = "휴우, 그래도 반이나 남았네."

# Row 429
$templateAL.Copy($ws.Range("A429:L429"))
$templateM.Copy($ws.Range("M429"))
$ws.Range("A429").Value = 45200.95795295139
$ws.Range("B429").Value = "molly7624@naver.com"
$ws.Range("C429").Value = "디지털미디어콘텐츠전공"
$ws.Range("D429").Value = 20211516
$ws.Range("E429").Value = "변재은"
$ws.Range("F429").Value = "민주 문자"
$ws.Range("G429").Value = "한글"
$ws.Range("H429").Value = "1개"
$ws.Range("I429").Value = 0.8
$ws.Range("J429").Value = "대한민국"
$ws.Range("K429").Value = "시간당 중위 임금이 60% 낮다"
$ws.Range("L429").Value = "Red"
$ws.Range("M429").Value = "휴우, 그래도 반이나 남았네."

# Row 430
$templateAL.Copy($ws.Range("A430:L430"))
$templateM.Copy($ws.Range("M430"))
$ws.Range("A430").Value = 45200.96560776621
$ws.Range("B430").Value = "choe0119@gmail.com"
$ws.Range("C430").Value = "의예과"
$ws.Range("D430").Value = 20226176
$ws.Range("E430").Value = "최태웅"
$ws.Range("F430").Value = "엘리트 문자"
$ws.Range("G430").Value = "한자"
$ws.Range("H430").Value = "1개"
$ws.Range("I430").Value = 0.2
$ws.Range("J430").Value = "영국"
$ws.Range("K430").Value = "건강이 좋지 않다"
$ws.Range("L430").Value = "Red"
$ws.Range("M430").Value = "헐, 반 밖에 안 남았네."

# Row 431
$templateAL.Copy($ws.Range("A431:L431"))
$templateN.Copy($ws.Range("N431"))
$ws.Range("A431").Value = 45200.970568564815
$ws.Range("B431").Value = "dksdksqh1018@naver.com"
$ws.Range("C431").Value = "미디어스쿨"
$ws.Range("D431").Value = 20232549
$ws.Range("E431").Value = "안보민"
$ws.Range("F431").Value = "민주 문자"
$ws.Range("G431").Value = "한글"
$ws.Range("H431").Value = "하나도 없다"
$ws.Range("I431").Value = 0.8
$ws.Range("J431").Value = "대한민국"
$ws.Range("K431").Value = "시간당 중위 임금이 60% 낮다"
$ws.Range("L431").Value = "Black"
$ws.Range("N431").Value = "헐, 반 밖에 안 남았네."

# Row 432
$templateAL.Copy($ws.Range("A432:L432"))
$templateN.Copy($ws.Range("N432"))
$ws.Range("A432").Value = 45200.97358670139
$ws.Range("B432").Value = "gaejisub@gmail.com"
$ws.Range("C432").Value = "콘텐츠it"
$ws.Range("D432").Value = 20225169
$ws.Range("E432").Value = "배승유"
$ws.Range("F432").Value = "민주 문자"
$ws.Range("G432").Value = "한글"
$ws.Range("H432").Value = "하나도 없다"
$ws.Range("I432").Value = 0.2
$ws.Range("J432").Value = "미국"
$ws.Range("K432").Value = "시간당 중위 임금이 60% 낮다"
$ws.Range("L432").Value = "Black"
$ws.Range("N432").Value = "휴우, 그래도 반이나 남았네."

# Row 433
$templateAL.Copy($ws.Range("A433:L433"))
$templateN.Copy($ws.Range("N433"))
$ws.Range("A433").Value = 45200.98353868056
$ws.Range("B433").Value = "gustj1654@naver.com"
$ws.Range("C433").Value = "심리학과"
$ws.Range("D433").Value = 20232137
$ws.Range("E433").Value = "조현서"
$ws.Range("F433").Value = "민주 문자"
$ws.Range("G433").Value = "한글"
$ws.Range("H433").Value = "1개"
$ws.Range("I433").Value = 0.8
$ws.Range("J433").Value = "대한민국"
$ws.Range("K433").Value = "시간당 중위 임금이 60% 낮다"
$ws.Range("L433").Value = "Black"
$ws.Range("N433").Value = "휴우, 그래도 반이나 남았네."

# Row 434
$templateAL.Copy($ws.Range("A434:L434"))
$templateM.Copy($ws.Range("M434"))
$ws.Range("A434").Value = 45200.98789751157
$ws.Range("B434").Value = "yejin4259@naver.com"
$ws.Range("C434").Value = "언어청각학부"
$ws.Range("D434").Value = 20233951
$ws.Range("E434").Value = "이예진"
$ws.Range("F434").Value = "민주 문자"
$ws.Range("G434").Value = "한글"
$ws.Range("H434").Value = "1개"
$ws.Range("I434").Value = 0.8
$ws.Range("J434").Value = "대한민국"
$ws.Range("K434").Value = "시간당 중위 임금이 60% 낮다"
$ws.Range("L434").Value = "Red"
$ws.Range("M434").Value = "휴우, 그래도 반이나 남았네."

# Row 435
$templateAL.Copy($ws.Range("A435:L435"))
$templateM.Copy($ws.Range("M435"))
$ws.Range("A435").Value = 45200.99348971064
$ws.Range("B435").Value = "rhy0787@naver.com"
$ws.Range("C435").Value = "식품영양학과"
$ws.Range("D435").Value = 20213827
$ws.Range("E435").Value = "유희영"
$ws.Range("F435").Value = "민주 문자"
$ws.Range("G435").Value = "한자"
$ws.Range("H435").Value = "하나도 없다"
$ws.Range("I435").Value = 0.8
$ws.Range("J435").Value = "대한민국"
$ws.Range("K435").Value = "시간당 중위 임금이 60% 낮다"
$ws.Range("L435").Value = "Red"
$ws.Range("M435").Value = "휴우, 그래도 반이나 남았네."

# Row 436
$templateAL.Copy($ws.Range("A436:L436"))
$templateM.Copy($ws.Range("M436"))
$ws.Range("A436").Value = 45200.99731702547
$ws.Range("B436").Value = "kxjenlee@naver.com"
$ws.Range("C436").Value = "글로벌비즈니스"
$ws.Range("D436").Value = 20226417
$ws.Range("E436").Value = "이제인"
$ws.Range("F436").Value = "민주 문자"
$ws.Range("G436").Value = "한글"
$ws.Range("H436").Value = "1개"
$ws.Range("I436").Value = 0.5
$ws.Range("J436").Value = "영국"
$ws.Range("K436").Value = "2배 정도 실직할 가능성이 높다"
$ws.Range("L436").Value = "Red"
$ws.Range("M436").Value = "모름/기타"

# Row 437
$templateAL.Copy($ws.Range("A437:L437"))
$templateM.Copy($ws.Range("M437"))
$ws.Range("A437").Value = 45201.00346453703
$ws.Range("B437").Value = "tjdus3641@gmail.com"
$ws.Range("C437").Value = "간호학과"
$ws.Range("D437").Value = 20226283
$ws.Range("E437").Value = "장서연"
$ws.Range("F437").Value = "민주 문자"
$ws.Range("G437").Value = "한글"
$ws.Range("H437").Value = "2개"
$ws.Range("I437").Value = 0.8
$ws.Range("J437").Value = "대한민국"
$ws.Range("K437").Value = "사회활동이나 자원활동에 덜 참여한다"
$ws.Range("L437").Value = "Red"
$ws.Range("M437").Value = "휴우, 그래도 반이나 남았네."

# Row 438
$templateAL.Copy($ws.Range("A438:L438"))
$templateN.Copy($ws.Range("N438"))
$ws.Range("A438").Value = 45201.00409956019
$ws.Range("B438").Value = "rkqls3333@gmail.com"
$ws.Range("C438").Value = "간호학과"
$ws.Range("D438").Value = 20236205
$ws.Range("E438").Value = "권가빈"
$ws.Range("F438").Value = "민주 문자"
$ws.Range("G438").Value = "한글"
$ws.Range("H438").Value = "2개"
$ws.Range("I438").Value = 0.2
$ws.Range("J438").Value = "대한민국"
$ws.Range("K438").Value = "2배 정도 실직할 가능성이 높다"
$ws.Range("L438").Value = "Black"
$ws.Range("N438").Value = "모름/기타"

